# Update "想去人数" (F column) counts on several rows across sheets to
# reflect the output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 3174
$ws1.Range("F7").Value = 2747
$ws1.Range("F12").Value = 286
$ws1.Range("F14").Value = 5696
$ws1.Range("F18").Value = 165
$ws1.Range("F21").Value = 1237
$ws1.Range("F24").Value = 124

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 1170
$ws2.Range("F6").Value = 242
$ws2.Range("F8").Value = 336
$ws2.Range("F13").Value = 634
$ws2.Range("F16").Value = 1000
$ws2.Range("F21").Value = 53
$ws2.Range("F23").Value = 340
$ws2.Range("F30").Value = 61
$ws2.Range("F33").Value = 41

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 2561
$ws3.Range("F6").Value = 1115
$ws3.Range("F9").Value = 1440

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2561
$ws4.Range("F6").Value = 1115
$ws4.Range("F7").Value = 1440
$ws4.Range("F14").Value = 3174
$ws4.Range("F15").Value = 2747
$ws4.Range("F19").Value = 242
$ws4.Range("F22").Value = 336
$ws4.Range("F25").Value = 5696
$ws4.Range("F30").Value = 634
$ws4.Range("F32").Value = 165
$ws4.Range("F38").Value = 53
$ws4.Range("F39").Value = 340
$ws4.Range("F40").Value = 1237
$ws4.Range("F46").Value = 41
